$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# Add Bill 11 row (row 12)
$ws.Range("A12").Value = 11
$ws.Range("C12").Value = "Burgstrips"
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 94
$ws.Range("F12").Value = "Akshay"
$ws.Range("H12").Value = 100
